# Swap the "daily" record data (Fecha/Calidad/Volumen/Precios/Unidad/Precio kg/Kg por unidad)
# between row pairs: (2,8), (4,7), (6,9) on the active sheet, leaving the
# identifying columns (A,B,C,E,F,G,H,I,J,K,R) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are swapped between the row pairs.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")

function Swap-Rows($ws, $rowA, $rowB, $cols) {
    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")
        $valA = $rangeA.Value2
        $valB = $rangeB.Value2
        $rangeA.Value2 = $valB
        $rangeB.Value2 = $valA
    }
}

Swap-Rows $ws 2 8 $cols
Swap-Rows $ws 4 7 $cols
Swap-Rows $ws 6 9 $cols
